# Trade #55 closed at 2026-02-17 12:52:10 - unknown UNKNOWN +0.000%
#
# Updates the Summary and Strategy Status roll-up figures to reflect the
# newly-closed trade, and appends the new trade row (#55 / row 56) to both
# the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.11   # Current Capital
$summary.Range("B4").Value = 0.1       # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 55        # Total Trades
$summary.Range("B8").Value = 19        # Losing Trades
$summary.Range("B9").Value = 41.82     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.11     # Capital
$status.Range("D4").Value = 55         # Trades
$status.Range("E4").Value = 0.1        # P&L $
$status.Range("F4").Value = 0.11       # P&L %
$status.Range("G4").Value = 41.82      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append new trade row (#55) to "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = 56

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$newRow").Value = 55

    # Date/time columns are stored as text, not Excel date/time serials -
    # force text number format before assignment so "2026-02-17" isn't
    # reinterpreted as a date.
    $ws.Range("B$newRow").NumberFormat = "@"
    $ws.Range("B$newRow").Value = "2026-02-17"
    $ws.Range("C$newRow").Value = "12:52:04"

    $ws.Range("D$newRow").Value = "MarketMaking"
    $ws.Range("E$newRow").Value = "DOWN"
    $ws.Range("F$newRow").Value = 0.83
    $ws.Range("G$newRow").Value = 0.75
    $ws.Range("H$newRow").Value = "CLOSED"
    $ws.Range("I$newRow").Value = -9.6386
    $ws.Range("J$newRow").Value = -0.08
    $ws.Range("K$newRow").Value = 100.11
    $ws.Range("L$newRow").Value = 0
    $ws.Range("M$newRow").Value = 0
    $ws.Range("N$newRow").Value = 0.6
    $ws.Range("O$newRow").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$newRow").Value = "early_exit"
    $ws.Range("Q$newRow").Value = 0.12
}
